$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 270 (shifts existing row 270+ down by one)
$ws.Rows.Item(270).Insert()

# Populate the newly inserted row 270 with the treatise datapacks title entry
$ws.Range("A270").Value = "settings.datapacks.title.treatise"
$ws.Range("B270").Value = "Treatise Datapacks"
